# Update the "想去人数" (wanted-to-go count) figures on both the
# "展览" and "全部类型" sheets to the newly scraped values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F14").Value = 279
    $ws.Range("F28").Value = 60
    $ws.Range("F35").Value = 247
    $ws.Range("F40").Value = 559
    $ws.Range("F42").Value = 3314
}
